# Adding ingest rake task, and version validations
#
# The template gains a new "version_publication_date" column and the two
# existing "edition_statement" / "publication_statement" headers are
# renamed to "version_edition_statement" / "version_publication_statement"
# (and given the "Comma" cell style, left-aligned) to make room for it.
# Everything to the right (additional_responsibility .. based_on_original)
# shifts one column over, from H:N to I:O.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new "version_publication_date" column by inserting a
# blank column at J (pushes additional_responsibility.. one column right).
$ws.Range("J1").EntireColumn.Insert()

# Rename the (still in place) H1/I1 headers to their "version_" variants,
# and fill in the brand-new J1 header.
$ws.Range("H1").Value = "version_edition_statement"
$ws.Range("I1").Value = "version_publication_statement"
$ws.Range("J1").Value = "version_publication_date"

# Apply the "Comma" cell style with left-aligned text to the two new
# "version_publication_*" headers (I1:J1).
$ws.Range("I1").Style = "Comma"
$ws.Range("I1").HorizontalAlignment = -4131
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
